# Add a new bulleted line ("How to release the memory of a dask object
# after .compute() is completed") as the very first item of the
# document's bullet list, immediately before the existing
# "In the sdoml_dataset notebook, continue debugging from where marked"
# paragraph.

$d = $word.ActiveDocument

# Find the paragraph that is currently first in the list so the new
# bullet can be inserted right before it, regardless of exactly where
# it sits in the body.
$searchRange = $d.Content
$searchRange.Find.Execute("In the sdoml_dataset notebook, continue debugging from where marked",
                           $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)

$targetPara = $searchRange.Paragraphs.First

# Inserting a paragraph mark before the target paragraph's range creates
# a brand-new paragraph that inherits the list/indent/font formatting
# (numPr ilvl=0/numId=1, 720/360 indent, 24pt) from the paragraph that
# follows it -- exactly the formatting the new bullet needs.
$targetPara.Range.InsertParagraphBefore()

# Fill in the text of the newly created (now first) paragraph.
$newPara = $d.Paragraphs.First
$newPara.Range.Text = "How to release the memory of a dask object after .compute() is completed"
